# HW 3 error fixes
# The "tuples" sheet had two stray SUM() helper columns (T and Y) left
# over from building the route table; remove them and let everything to
# their right slide over one column (Y->X, X->W, W->V, V->U).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column T only contained =SUM(...) helper cells (rows 5,8,10,12,15) -
# deleting it shifts U:Y left to T:X.
$ws.Columns("T").Delete()

# What was column Y (the grand-total =SUM(...) cells) is now column X
# after the shift above - delete it too, shifting the route-table data
# (originally V:X) left into U:W.
$ws.Columns("X").Delete()

# Leave the selection on the now-empty column X, matching the saved
# workbook state.
$ws.Columns("X").Select() | Out-Null
